# Add a new review row (row 8) to the reviewdb sheet, mirroring the layout
# of the existing rows (appid/keyword/email/recovery-email/time/review),
# including the two mailto: hyperlinks on the email columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 8 with row 7's formatting (column styles: A/F=1, B/E=0, C/D=2)
# so the new cells pick up the same cellXfs indices as the rest of the table.
$ws.Range("A7:F7").Copy($ws.Range("A8:F8"))
$excel.CutCopyMode = $false

# Populate the new row's values.
$ws.Range("A8").Value = "com.singleton.strechy"
$ws.Range("B8").Value = "stretchy"
$ws.Range("C8").Value = "nitanfriman@gmail.com"
$ws.Range("D8").Value = "ronoren61@gmail.com"
$ws.Range("E8").Value = "27/5/2019 15:59"
$ws.Range("F8").Value = "Hilarious to play this game when you are drunk. Try it some time haha"

# Wire up the mailto hyperlinks for the two email columns, same as every
# other row in the sheet.
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:nitanfriman@gmail.com", [Type]::Missing, [Type]::Missing, "nitanfriman@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:ronoren61@gmail.com", [Type]::Missing, [Type]::Missing, "ronoren61@gmail.com")

# Hyperlinks.Add() re-styles the target cell with the built-in "Hyperlink"
# style; restore the table's normal per-column formatting afterwards so the
# cells keep using the existing style indices instead of a new one.
$ws.Range("C7").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection onto the newly added last cell, matching what
# Excel does after typing data into a new row.
[void]$ws.Range("F8").Select()
